$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 460
$ws.Range("I33").Value = 566.6667
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 566.6667
$ws.Range("L33").Value = 300
$ws.Range("M33").Value = -337.6667
$ws.Range("N33").Value = -758
$ws.Range("H107").Value = 568.1786
$ws.Range("I107").Value = 515.7308
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 515.7308
$ws.Range("L107").Value = 1250
$ws.Range("M107").Value = 1404.2692
$ws.Range("N107").Value = -5090
$ws.Range("H137").Value = 557879.25
$ws.Range("I137").Value = 1483.3103
$ws.Range("J137").Value = 1326235.5
$ws.Range("K137").Value = 4449.9309
$ws.Range("L137").Value = 3978706.5
$ws.Range("M137").Value = -1899.9309
$ws.Range("N137").Value = -3983806.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4617.0566
$ws.Range("I32").Value = 3646.4404
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 3646.4404
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -3359.4404
$ws.Range("N32").Value = -25574
$ws.Range("H61").Value = 5157.6772
$ws.Range("I61").Value = 3858.3333
$ws.Range("J61").Value = 13928.25
$ws.Range("K61").Value = 3858.3333
$ws.Range("L61").Value = 13928.25
$ws.Range("M61").Value = -3646.3333
$ws.Range("N61").Value = -14352.25
$ws.Range("H74").Value = 4960.552
$ws.Range("I74").Value = 1952.6316
$ws.Range("J74").Value = 10675.6
$ws.Range("K74").Value = 1952.6316
$ws.Range("L74").Value = 10675.6
$ws.Range("M74").Value = -1078.6316
$ws.Range("N74").Value = -12423.6
$ws.Range("H77").Value = 4960.552
$ws.Range("I77").Value = 1952.6316
$ws.Range("J77").Value = 10675.6
$ws.Range("K77").Value = 9763.157999999999
$ws.Range("L77").Value = 53378
$ws.Range("M77").Value = -5395.157999999999
$ws.Range("N77").Value = -62114
$ws.Range("H122").Value = 5210445
$ws.Range("I122").Value = 2685.0908
$ws.Range("J122").Value = 9617011
$ws.Range("K122").Value = 8055.2724
$ws.Range("L122").Value = 28851033
$ws.Range("M122").Value = -5605.2724
$ws.Range("N122").Value = -28855933
$ws.Range("H132").Value = 3366.9104
$ws.Range("I132").Value = 1029.62
$ws.Range("K132").Value = 3088.86
$ws.Range("M132").Value = -558.8599999999997
$ws.Range("H136").Value = 5157.6772
$ws.Range("I136").Value = 3858.3333
$ws.Range("J136").Value = 13928.25
$ws.Range("K136").Value = 11574.9999
$ws.Range("L136").Value = 41784.75
$ws.Range("M136").Value = -9024.999899999999
$ws.Range("N136").Value = -46884.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 954845.4
$ws.Range("I105").Value = 1425054.5
$ws.Range("J105").Value = 14427.182
$ws.Range("K105").Value = 1425054.5
$ws.Range("L105").Value = 14427.182
$ws.Range("M105").Value = -1423307.5
$ws.Range("N105").Value = -17921.182
$ws.Range("H134").Value = 6650.237
$ws.Range("I134").Value = 6878.893
$ws.Range("J134").Value = 6010
$ws.Range("K134").Value = 20636.679
$ws.Range("L134").Value = 18030
$ws.Range("M134").Value = -18101.679
$ws.Range("N134").Value = -23100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2091.926
$ws.Range("I31").Value = 1426.7441
$ws.Range("K31").Value = 1426.7441
$ws.Range("M31").Value = -1131.7441
$ws.Range("H34").Value = 2091.926
$ws.Range("I34").Value = 1426.7441
$ws.Range("K34").Value = 1426.7441
$ws.Range("M34").Value = -1224.7441
$ws.Range("H86").Value = 2068
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2068
$ws.Range("K86").Value = 0
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = 2068
$ws.Range("N86").Value = -4314
$ws.Range("H89").Value = 2068
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2068
$ws.Range("K89").Value = 0
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = 10340
$ws.Range("N89").Value = -21572
$ws.Range("H132").Value = 2699.5186
$ws.Range("I132").Value = 2010.2941
$ws.Range("J132").Value = 3871.2
$ws.Range("K132").Value = 6030.8823
$ws.Range("L132").Value = 11613.6
$ws.Range("M132").Value = -3500.8823
$ws.Range("N132").Value = -16673.6
$ws.Range("H134").Value = 3259.8723
$ws.Range("I134").Value = 2056.1904
$ws.Range("J134").Value = 4232.077
$ws.Range("K134").Value = 6168.5712
$ws.Range("L134").Value = 12696.231
$ws.Range("M134").Value = -3633.5712
$ws.Range("N134").Value = -17766.231

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 3327.0908
$ws.Range("I75").Value = 1797
$ws.Range("J75").Value = 3480.1
$ws.Range("K75").Value = 5391
$ws.Range("L75").Value = 10440.3
$ws.Range("M75").Value = -4393
$ws.Range("N75").Value = -12436.3
$ws.Range("H78").Value = 3327.0908
$ws.Range("I78").Value = 1797
$ws.Range("J78").Value = 3480.1
$ws.Range("K78").Value = 16173
$ws.Range("L78").Value = 31320.9
$ws.Range("M78").Value = -11181
$ws.Range("N78").Value = -41304.89999999999
$ws.Range("H108").Value = 3955.111
$ws.Range("I108").Value = 1508
$ws.Range("J108").Value = 7014
$ws.Range("K108").Value = 4524
$ws.Range("L108").Value = 21042
$ws.Range("M108").Value = -1644
$ws.Range("N108").Value = -26802
$ws.Range("H120").Value = 7237.55
$ws.Range("I120").Value = 7533.3335
$ws.Range("J120").Value = 7185.353
$ws.Range("K120").Value = 22600.0005
$ws.Range("L120").Value = 21556.059
$ws.Range("M120").Value = -17762.0005
$ws.Range("N120").Value = -31232.059

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2251.647
$ws.Range("I132").Value = 2252.4167
$ws.Range("K132").Value = 6757.250100000001
$ws.Range("M132").Value = -4227.250100000001
$ws.Range("H134").Value = 39580.832
$ws.Range("J134").Value = 39580.832
$ws.Range("L134").Value = 118742.496
$ws.Range("N134").Value = -123812.496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 940.7692
$ws.Range("J22").Value = 916.36365
$ws.Range("L22").Value = 916.36365
$ws.Range("N22").Value = -1506.36365
$ws.Range("H27").Value = 940.7692
$ws.Range("J27").Value = 916.36365
$ws.Range("L27").Value = 916.36365
$ws.Range("N27").Value = -1130.36365
$ws.Range("H46").Value = 699.8889
$ws.Range("I46").Value = 699.6667
$ws.Range("J46").Value = 700
$ws.Range("K46").Value = 699.6667
$ws.Range("L46").Value = 700
$ws.Range("M46").Value = -511.6667
$ws.Range("N46").Value = -1076
$ws.Range("H55").Value = 133697.56
$ws.Range("I55").Value = 210876.89
$ws.Range("J55").Value = 387.81818
$ws.Range("K55").Value = 210876.89
$ws.Range("L55").Value = 387.81818
$ws.Range("M55").Value = -210703.89
$ws.Range("N55").Value = -733.81818
$ws.Range("H82").Value = 1595.8572
$ws.Range("J82").Value = 2070
$ws.Range("L82").Value = 2070
$ws.Range("N82").Value = -2792
$ws.Range("H85").Value = 1595.8572
$ws.Range("J85").Value = 2070
$ws.Range("L85").Value = 2070
$ws.Range("N85").Value = -4566
$ws.Range("H93").Value = 874.875
$ws.Range("I93").Value = 939.8
$ws.Range("J93").Value = 766.6667
$ws.Range("K93").Value = 939.8
$ws.Range("L93").Value = 766.6667
$ws.Range("M93").Value = 308.2
$ws.Range("N93").Value = -3262.6667
$ws.Range("H122").Value = 4007.475
$ws.Range("I122").Value = 3882.3794
$ws.Range("J122").Value = 4337.273
$ws.Range("K122").Value = 11647.1382
$ws.Range("L122").Value = 13011.819
$ws.Range("M122").Value = -9197.138199999999
$ws.Range("N122").Value = -17911.819
$ws.Range("H132").Value = 4226.222
$ws.Range("I132").Value = 3851.6667
$ws.Range("J132").Value = 4750.6
$ws.Range("K132").Value = 11555.0001
$ws.Range("L132").Value = 14251.8
$ws.Range("M132").Value = -9025.000100000001
$ws.Range("N132").Value = -19311.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1390.7106
$ws.Range("I132").Value = 1583.25
$ws.Range("J132").Value = 1060.6428
$ws.Range("K132").Value = 4749.75
$ws.Range("L132").Value = 3181.9284
$ws.Range("M132").Value = -2219.75
$ws.Range("N132").Value = -8241.928400000001
$ws.Range("H136").Value = 3975.4307
$ws.Range("I136").Value = 1987.8
$ws.Range("J136").Value = 6294.3335
$ws.Range("K136").Value = 5963.4
$ws.Range("L136").Value = 18883.0005
$ws.Range("M136").Value = -3413.4
$ws.Range("N136").Value = -23983.0005
